$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "平潭发展"
$ws.Range("B2").Value = "平潭发展"

$ws.Range("A3").Value = "孚日股份"
$ws.Range("C3").Value = "安泰集团"

$ws.Range("A4").Value = "多氟多"
$ws.Range("B4").Value = "长城军工"
$ws.Range("C4").Value = "孚日股份"

$ws.Range("A5").Value = "华胜天成"
$ws.Range("B5").Value = "安泰集团"
$ws.Range("C5").Value = "雪人集团"

$ws.Range("A6").Value = "安泰集团"
$ws.Range("B6").Value = "三六零"
$ws.Range("C6").Value = "多氟多"

$ws.Range("A7").Value = "长城军工"
$ws.Range("B7").Value = "华胜天成"
$ws.Range("C7").Value = "海马汽车"

$ws.Range("A8").Value = "海峡创新"
$ws.Range("B8").Value = "航天发展"
$ws.Range("C8").Value = "华胜天成"

$ws.Range("A9").Value = "永太科技"
$ws.Range("B9").Value = "华夏幸福"
$ws.Range("C9").Value = "海峡创新"

$ws.Range("A10").Value = "航天发展"
$ws.Range("B10").Value = "日出东方"
$ws.Range("C10").Value = "三木集团"

$ws.Range("A11").Value = "雪人集团"
$ws.Range("B11").Value = "多氟多"
$ws.Range("C11").Value = "长城军工"

$ws.Range("A12").Value = "日出东方"
$ws.Range("B12").Value = "雪人集团"
$ws.Range("C12").Value = "永太科技"

$ws.Range("A13").Value = "三六零"
$ws.Range("B13").Value = "天齐锂业"
$ws.Range("C13").Value = "华夏幸福"

$ws.Range("A14").Value = "华夏幸福"
$ws.Range("B14").Value = "海峡创新"
$ws.Range("C14").Value = "三六零"

$ws.Range("B15").Value = "盛新锂能"
$ws.Range("C15").Value = "盈新发展"

$ws.Range("A16").Value = "天际股份"
$ws.Range("B16").Value = "盈新发展"
$ws.Range("C16").Value = "摩恩电气"

$ws.Range("A17").Value = "盈新发展"
$ws.Range("B17").Value = "海马汽车"
$ws.Range("C17").Value = "东百集团"

$ws.Range("A18").Value = "天齐锂业"
$ws.Range("B18").Value = "众生药业"
$ws.Range("C18").Value = "天际股份"

$ws.Range("A19").Value = "三木集团"
$ws.Range("B19").Value = "永太科技"
$ws.Range("C19").Value = "大有能源"

$ws.Range("B20").Value = "大有能源"
$ws.Range("C20").Value = "航天发展"

$ws.Range("A21").Value = "盛新锂能"
$ws.Range("C21").Value = "人民同泰"
